$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (reflects "Through 2021-12-16")
$ws.Name = "Through 2021-12-16"

# Row 10 (August) updates
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 74
$ws.Range("G10").Value = 0.0633

# Row 14 (December through 12-16) updates
$ws.Range("A14").Value = "December (through 12-16)"
$ws.Range("C14").Value = 16
$ws.Range("D14").Value = 0.1579
$ws.Range("F14").Value = 45
$ws.Range("G14").Value = 0.1
$ws.Range("I14").Value = 52
$ws.Range("J14").Value = 0.0877
$ws.Range("L14").Value = 32
$ws.Range("M14").Value = 0.0857
$ws.Range("O14").Value = 24
$ws.Range("P14").Value = 0.1111
$ws.Range("R14").Value = 74
$ws.Range("S14").Value = 0.0513
$ws.Range("U14").Value = 124
$ws.Range("V14").Value = 0.008

# Row 15 (Total) updates
$ws.Range("C15").Value = 274
$ws.Range("D15").Value = 0.1161
$ws.Range("E15").Value = 64
$ws.Range("F15").Value = 549
$ws.Range("G15").Value = 0.1044
$ws.Range("I15").Value = 810
$ws.Range("J15").Value = 0.0774
$ws.Range("L15").Value = 640
$ws.Range("M15").Value = 0.1074
$ws.Range("O15").Value = 504
$ws.Range("P15").Value = 0.1016
$ws.Range("R15").Value = 1274
$ws.Range("U15").Value = 1666
$ws.Range("V15").Value = 0.0577
